# Insert two new data rows (for the new market date 44559) right before the
# existing row 263, shifting all subsequent rows (old 263..348) down to
# 265..350. Then populate the two freshly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 263 (each Insert() pushes rows 263+ down by one)
$ws.Rows.Item(263).Insert()
$ws.Rows.Item(263).Insert()

# New row 263: Betarraga, "Primera" quality, market date 44559 (2021-12-29)
$ws.Cells.Item(263, 1).Value = 9
$ws.Cells.Item(263, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(263, 3).Value = "Metropolitana"
$ws.Cells.Item(263, 4).Value = 44559
$ws.Cells.Item(263, 5).Value = 13
$ws.Cells.Item(263, 6).Value = 100114014
$ws.Cells.Item(263, 7).Value = "Betarraga"
$ws.Cells.Item(263, 8).Value = "Sin especificar"
$ws.Cells.Item(263, 9).Value = "Primera"
$ws.Cells.Item(263, 10).Value = 6100
$ws.Cells.Item(263, 11).Value = 80
$ws.Cells.Item(263, 12).Value = 90
$ws.Cells.Item(263, 13).Value = 85
$ws.Cells.Item(263, 14).Value = "$/unidad"
$ws.Cells.Item(263, 15).Value = "Región Metropolitana"
$ws.Cells.Item(263, 16).Value = 85
$ws.Cells.Item(263, 17).Value = 1
$ws.Cells.Item(263, 18).Value = "Hortaliza"

# New row 264: Betarraga, "Segunda" quality, market date 44559 (2021-12-29)
$ws.Cells.Item(264, 1).Value = 9
$ws.Cells.Item(264, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(264, 3).Value = "Metropolitana"
$ws.Cells.Item(264, 4).Value = 44559
$ws.Cells.Item(264, 5).Value = 13
$ws.Cells.Item(264, 6).Value = 100114014
$ws.Cells.Item(264, 7).Value = "Betarraga"
$ws.Cells.Item(264, 8).Value = "Sin especificar"
$ws.Cells.Item(264, 9).Value = "Segunda"
$ws.Cells.Item(264, 10).Value = 2500
$ws.Cells.Item(264, 11).Value = 60
$ws.Cells.Item(264, 12).Value = 70
$ws.Cells.Item(264, 13).Value = 65
$ws.Cells.Item(264, 14).Value = "$/unidad"
$ws.Cells.Item(264, 15).Value = "Región Metropolitana"
$ws.Cells.Item(264, 16).Value = 65
$ws.Cells.Item(264, 17).Value = 1
$ws.Cells.Item(264, 18).Value = "Hortaliza"
